$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the logical-operator formulas for rows 2-30 in columns G, H, I, J
for ($r = 2; $r -le 30; $r++) {
    $ws.Range("G$r").Formula = "=IF(F$r>=75000,""big"",""small"")"
    $ws.Range("H$r").Formula = "=IF(OR(C$r=""MacBook Pro"", C$r=""MacBook Air""),""Macbook"", IF(C$r=""Mac Pro"", ""Mac"", """"))"
    $ws.Range("I$r").Formula = "=IF(AND(C$r=""MacBook Pro"", F$r>=35000), TRUE, """")"
    $ws.Range("J$r").Formula = "=IF(I$r=TRUE, F$r, """")"
}

# Summary formulas
$ws.Range("M16").Formula = "=COUNTIF(I2:I30,TRUE)"
$ws.Range("M17").Formula = "=SUMIF(I2:I30, TRUE, J2:J30)"
$ws.Range("M18").Formula = "=COUNTIF(G2:G30,""small"")"

$wb.Save()
